$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from E1 (bold, bordered, centered) into the new F1 header cell,
# then overwrite its value with the new header text "Modelo".
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Modelo"

# Add the new data cell F2 with the model pipeline description (no special style, like E2).
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
